$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Update stack-trace line numbers to match the 3.0.0 -> 3.1.0 version bump.

Replace-Text "PaginationServices.asStyle(PaginationServices.java:155)" "PaginationServices.asStyle(PaginationServices.java:207)"

Replace-Text "M2DocEvaluator.caseQuery(M2DocEvaluator.java:559)" "M2DocEvaluator.caseQuery(M2DocEvaluator.java:586)"

# All three occurrences of doSwitch(M2DocEvaluator.java:1216) become 1239.
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1216)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1239)"

Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1425)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1464)"

Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:287)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:296)"

Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:276)" "M2DocEvaluator.generate(M2DocEvaluator.java:281)"

Replace-Text "M2DocUtils.generate(M2DocUtils.java:694)" "M2DocUtils.generate(M2DocUtils.java:805)"

Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)"

Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:389)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:420)"

# Insert a new stack frame (RunBefores.evaluate) right before the second
# occurrence of "RunAfters.evaluate(RunAfters.java:27)" (the one that is
# immediately preceded by "ParentRunner$2.evaluate(ParentRunner.java:268)"
# and followed by "ParentRunner.run(ParentRunner.java:363)" then
# "Suite.runChild(Suite.java:128)").
$find = "ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)"
$replace = "ParentRunner`$2.evaluate(ParentRunner.java:268)`n`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)`n`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)`n`tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)`n`tat org.junit.runners.Suite.runChild(Suite.java:128)"
Replace-Text $find $replace
